# Apply cryptos list price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B-E hold text (coin name / link / price / volume%) even when the
# price text looks numeric (e.g. "248.04"), so force text formatting while
# writing, then drop the number format again to leave cell styling untouched.
$textCols = $ws.Range("B2:E51")
$textCols.NumberFormat = "@"

$ws.Range("D2").Value = "37.123.92"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "2.050.30"
$ws.Range("E3").Value = "  -1.25%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "248.04"
$ws.Range("E5").Value = "  -1.95%  "

$ws.Range("D6").Value = "0.663"
$ws.Range("E6").Value = "  -1.61%  "

$ws.Range("D7").Value = "57.86"
$ws.Range("E7").Value = "  -1.97%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -2.29%  "

$ws.Range("D10").Value = "0.0774"
$ws.Range("E10").Value = "  -2.07%  "

$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("D12").Value = "15.90"
$ws.Range("E12").Value = "  -0.85%  "

$ws.Range("D13").Value = "0.871"
$ws.Range("E13").Value = "  +6.75%  "

$ws.Range("D14").Value = "2.348.94"
$ws.Range("E14").Value = "  -1.34%  "

$ws.Range("D15").Value = "5.70"
$ws.Range("E15").Value = "  +2.69%  "

$ws.Range("D16").Value = "2.051.88"
$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("D17").Value = "17.92"
$ws.Range("E17").Value = "  +15.31%  "

$ws.Range("D18").Value = "37.130.03"
$ws.Range("E18").Value = "  +0.14%  "

$ws.Range("D19").Value = "74.90"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").Value = "0.0₃0891"
$ws.Range("E20").Value = "  -3.78%  "

$ws.Range("D21").Value = "5.36"
$ws.Range("E21").Value = "  -1.62%  "

$ws.Range("D22").Value = "236.91"
$ws.Range("E22").Value = "  -0.95%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("E24").Value = "  +1.98%  "

$ws.Range("D25").Value = "9.51"
$ws.Range("E25").Value = "  +2.22%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "2.17"
$ws.Range("E26").Value = "  -5.25%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "169.16"
$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("D28").Value = "20.05"
$ws.Range("E28").Value = "  -1.41%  "

$ws.Range("E29").Value = "  -1.42%  "

$ws.Range("D30").Value = "4.80"
$ws.Range("E30").Value = "  -0.24%  "

$ws.Range("D31").Value = "1.12"
$ws.Range("E31").Value = "  -0.92%  "

$ws.Range("D32").Value = "0.0617"
$ws.Range("E32").Value = "  -2.90%  "

$ws.Range("D33").Value = "4.47"
$ws.Range("E33").Value = "  +0.51%  "

$ws.Range("D34").Value = "0.0895"
$ws.Range("E34").Value = "  -2.26%  "

$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("D36").Value = "2.25"
$ws.Range("E36").Value = "  -2.05%  "

$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").Value = "3.26"
$ws.Range("E38").Value = "  +15.52%  "

$ws.Range("D39").Value = "1.34"
$ws.Range("E39").Value = "  -2.12%  "

$ws.Range("D40").Value = "5.16"
$ws.Range("E40").Value = "  +15.35%  "

$ws.Range("E41").Value = "  -2.09%  "

$ws.Range("D42").Value = "17.27"
$ws.Range("E42").Value = "  -4.05%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "1.14"
$ws.Range("E43").Value = "  -2.46%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.0943"
$ws.Range("E44").Value = "  -20.96%  "

$ws.Range("D45").Value = "95.93"
$ws.Range("E45").Value = "  -2.72%  "

$ws.Range("D46").Value = "2.43"
$ws.Range("E46").Value = "  -1.67%  "

$ws.Range("D47").Value = "1.275.34"
$ws.Range("E47").Value = "  -2.12%  "

$ws.Range("D48").Value = "2.86"
$ws.Range("E48").Value = "  -3.59%  "

$ws.Range("D49").Value = "6.81"
$ws.Range("E49").Value = "  -2.11%  "

$ws.Range("D50").Value = "2.233.22"
$ws.Range("E50").Value = "  -1.29%  "

$ws.Range("D51").Value = "43.63"
$ws.Range("E51").Value = "  -0.36%  "

$textCols.ClearFormats()